$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 50 new patients appended below the existing 150 (rows 152-201),
# matching the row height used by the rest of the list.
$ws.Range("A152:B201").RowHeight = 13.8

# Column A holds numeric-looking ids that must stay text (shared string),
# like the rest of the id column, so force Text format while writing them.
$ws.Range("A152:A201").NumberFormat = "@"

$ws.Range("A152").Value = '13667'
$ws.Range("B152").Value = 'GARRISON,Janet'
$ws.Range("A153").Value = '13702'
$ws.Range("B153").Value = 'STEIN,Janet'
$ws.Range("A154").Value = '13737'
$ws.Range("B154").Value = 'PARSONS,Alice'
$ws.Range("A155").Value = '13773'
$ws.Range("B155").Value = 'BASS,Megan'
$ws.Range("A156").Value = '13809'
$ws.Range("B156").Value = 'GARDNER,Stephen'
$ws.Range("A157").Value = '13844'
$ws.Range("B157").Value = 'MCCLAIN,Hannah'
$ws.Range("A158").Value = '13879'
$ws.Range("B158").Value = 'BULLOCK,Debra'
$ws.Range("A159").Value = '13914'
$ws.Range("B159").Value = 'AVILA,Brittany'
$ws.Range("A160").Value = '13949'
$ws.Range("B160").Value = 'BALL,Nicholas'
$ws.Range("A161").Value = '13985'
$ws.Range("B161").Value = 'REYNOLDS,Kelly'
$ws.Range("A162").Value = '14020'
$ws.Range("B162").Value = 'BENDER,Amy'
$ws.Range("A163").Value = '14056'
$ws.Range("B163").Value = 'WATKINS,Rose'
$ws.Range("A164").Value = '14091'
$ws.Range("B164").Value = 'DILLON,Betty'
$ws.Range("A165").Value = '14126'
$ws.Range("B165").Value = 'CHANG,Betty'
$ws.Range("A166").Value = '14162'
$ws.Range("B166").Value = 'FIGUEROA,Debra'
$ws.Range("A167").Value = '14198'
$ws.Range("B167").Value = 'FLETCHER,Rachel'
$ws.Range("A168").Value = '14233'
$ws.Range("B168").Value = 'MCCULLOUGH,Heather'
$ws.Range("A169").Value = '14269'
$ws.Range("B169").Value = 'SHERMAN,Evelyn'
$ws.Range("A170").Value = '14304'
$ws.Range("B170").Value = 'WARREN,Nicole'
$ws.Range("A171").Value = '14340'
$ws.Range("B171").Value = 'LEON,Andrea'
$ws.Range("A172").Value = '14375'
$ws.Range("B172").Value = 'KLINE,Lauren'
$ws.Range("A173").Value = '14411'
$ws.Range("B173").Value = 'COCHRAN,Virginia'
$ws.Range("A174").Value = '14447'
$ws.Range("B174").Value = 'HENSLEY,Beverly'
$ws.Range("A175").Value = '14482'
$ws.Range("B175").Value = 'HATFIELD,Anna'
$ws.Range("A176").Value = '14518'
$ws.Range("B176").Value = 'MEYERS,Marilyn'
$ws.Range("A177").Value = '14554'
$ws.Range("B177").Value = 'RANDALL,Amber'
$ws.Range("A178").Value = '14589'
$ws.Range("B178").Value = 'FINLEY,Olivia'
$ws.Range("A179").Value = '14625'
$ws.Range("B179").Value = 'COLON,Laura'
$ws.Range("A180").Value = '14660'
$ws.Range("B180").Value = 'MUELLER,Dorothy'
$ws.Range("A181").Value = '14696'
$ws.Range("B181").Value = 'CANNON,Megan'
$ws.Range("A182").Value = '14731'
$ws.Range("B182").Value = 'BENTON,Marie'
$ws.Range("A183").Value = '14766'
$ws.Range("B183").Value = 'WARREN,Olivia'
$ws.Range("A184").Value = '14801'
$ws.Range("B184").Value = 'BISHOP,Kathryn'
$ws.Range("A185").Value = '14837'
$ws.Range("B185").Value = 'FISCHER,Jean'
$ws.Range("A186").Value = '14872'
$ws.Range("B186").Value = 'SWEENEY,Amanda'
$ws.Range("A187").Value = '14908'
$ws.Range("B187").Value = 'BARTLETT,Dorothy'
$ws.Range("A188").Value = '14944'
$ws.Range("B188").Value = 'CLARKE,Martha'
$ws.Range("A189").Value = '14979'
$ws.Range("B189").Value = 'TUCKER,Linda'
$ws.Range("A190").Value = '15014'
$ws.Range("B190").Value = 'OLSEN,Sandra'
$ws.Range("A191").Value = '15049'
$ws.Range("B191").Value = 'SAMPSON,Sara'
$ws.Range("A192").Value = '15085'
$ws.Range("B192").Value = 'FLYNN,Janet'
$ws.Range("A193").Value = '15121'
$ws.Range("B193").Value = 'BLAKE,Melissa'
$ws.Range("A194").Value = '15157'
$ws.Range("B194").Value = 'CLEMENTS,Stephen'
$ws.Range("A195").Value = '15193'
$ws.Range("B195").Value = 'CURTIS,Victoria'
$ws.Range("A196").Value = '15229'
$ws.Range("B196").Value = 'MOORE,Sarah'
$ws.Range("A197").Value = '15264'
$ws.Range("B197").Value = 'RANDOLPH,Andrea'
$ws.Range("A198").Value = '15300'
$ws.Range("B198").Value = 'PETTY,Cheryl'
$ws.Range("A199").Value = '15335'
$ws.Range("B199").Value = 'FISCHER,Carol'
$ws.Range("A200").Value = '15370'
$ws.Range("B200").Value = 'MITCHELL,Victoria'
$ws.Range("A201").Value = '15405'
$ws.Range("B201").Value = 'BRADFORD,Steven'

# Restore the default General format now that the text values are committed
$ws.Range("A152:A201").NumberFormat = "General"

$null = $ws.Range("A152:B201").Select()
